$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look numeric must stay stored as text
# (matching the source data's inline-string cells), so set NumberFormat
# to Text ("@") before writing them - otherwise Excel would silently convert
# strings like "88.00" into the number 88 and drop the trailing zero.
$textCells = @('D5','D6','D8','D9','D10','D11','D12','D13','D14','D19','D20','D21','D22','D23','D24','D25','D26','D27','D28','D29','D30','D31','D32','D34','D36','D38','D39','D40','D41','D46','D48','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '69.633.46'
$ws.Range('E2').Value = '  +4.97%  '
$ws.Range('D3').Value = '3.635.17'
$ws.Range('E3').Value = '  +18.17%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '590.74'
$ws.Range('E5').Value = '  +2.68%  '
$ws.Range('D6').Value = '186.44'
$ws.Range('E6').Value = '  +8.99%  '
$ws.Range('D7').Value = '3.631.77'
$ws.Range('E7').Value = '  +18.18%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  +4.79%  '
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  +9.64%  '
$ws.Range('D11').Value = '6.56'
$ws.Range('E11').Value = '  +5.15%  '
$ws.Range('D12').Value = '0.498'
$ws.Range('E12').Value = '  +5.96%  '
$ws.Range('D13').Value = '39.41'
$ws.Range('E13').Value = '  +10.12%  '
$ws.Range('D14').Value = '0.0000255'
$ws.Range('E14').Value = '  +6.75%  '
$ws.Range('D15').Value = '4.243.49'
$ws.Range('E15').Value = '  +18.26%  '
$ws.Range('D16').Value = '3.638.80'
$ws.Range('E16').Value = '  +18.30%  '
$ws.Range('D17').Value = '69.728.08'
$ws.Range('E17').Value = '  +5.14%  '
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').Value = '7.54'
$ws.Range('E19').Value = '  +8.40%  '
$ws.Range('D20').Value = '17.24'
$ws.Range('E20').Value = '  +3.59%  '
$ws.Range('D21').Value = '509.71'
$ws.Range('E21').Value = '  +4.66%  '
$ws.Range('D22').Value = '9.29'
$ws.Range('E22').Value = '  +20.98%  '
$ws.Range('D23').Value = '0.750'
$ws.Range('E23').Value = '  +9.27%  '
$ws.Range('D24').Value = '88.00'
$ws.Range('E24').Value = '  +6.77%  '
$ws.Range('D25').Value = '13.59'
$ws.Range('E25').Value = '  +7.58%  '
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  +9.42%  '
$ws.Range('D27').Value = '10.88'
$ws.Range('E27').Value = '  +7.63%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '2.56'
$ws.Range('E29').Value = '  +13.61%  '
$ws.Range('D30').Value = '8.26'
$ws.Range('E30').Value = '  +5.23%  '
$ws.Range('D31').Value = '32.28'
$ws.Range('E31').Value = '  +16.43%  '
$ws.Range('D32').Value = '2.75'
$ws.Range('E32').Value = '  +5.80%  '
$ws.Range('E33').Value = '  +18.99%  '
$ws.Range('D34').Value = '0.118'
$ws.Range('E34').Value = '  +5.89%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '6.19'
$ws.Range('E36').Value = '  +11.42%  '
$ws.Range('E37').Value = '  +8.31%  '
$ws.Range('D38').Value = '0.336'
$ws.Range('E38').Value = '  +11.32%  '
$ws.Range('B39').Value = 'Arweave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D39').Value = '47.11'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '2.10'
$ws.Range('E40').Value = '  +7.49%  '
$ws.Range('D41').Value = '50.75'
$ws.Range('E41').Value = '  +3.54%  '
$ws.Range('E42').Value = '  +4.70%  '
$ws.Range('D43').Value = '3.151.35'
$ws.Range('E43').Value = '  +13.71%  '
$ws.Range('E44').Value = '  +7.27%  '
$ws.Range('E45').Value = '  +10.78%  '
$ws.Range('D46').Value = '405.06'
$ws.Range('E46').Value = '  +10.91%  '
$ws.Range('E47').Value = '  +6.69%  '
$ws.Range('D48').Value = '27.96'
$ws.Range('E48').Value = '  +15.44%  '
$ws.Range('D49').Value = '136.64'
$ws.Range('E49').Value = '  +1.40%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '2.46'
$ws.Range('E50').Value = '  +14.55%  '
$ws.Range('B51').Value = 'USDe'
$ws.Range('C51').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.03%  '
